# Weekly epidemiological-event data refresh ("envio semana 20 de 2025"):
# updated Esperado/Observado/valor p counts for each evento, plus one
# brand-new evento row (205 - Chagas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Chagas" (evento 205) is a newly reported event this week; insert its row
# before the existing "Dengue" (210) row so the table stays sorted by
# evento code, shifting every following row down by one.
$ws.Rows.Item(5).Insert()

# Column A holds evento codes as text (e.g. "113"); force text formatting
# on the new row so the numeric-looking code "205" is not auto-converted
# to a number.
$ws.Range("A5").NumberFormat = "@"

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.03

$ws.Range("A5").Value = "205"
$ws.Range("B5").Value = "Chagas"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 10

$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.06

$ws.Range("C9").Value = 42
$ws.Range("D9").Value = 51
$ws.Range("E9").Value = 0.02

$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 0.06

$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 0

$ws.Range("C13").Value = 7

$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0.37

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 0

$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 0.11

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0.27

$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0.16

$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0.06

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1

$ws.Range("C23").Value = 7
$ws.Range("E23").Value = 0.13

$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0.14

$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1

$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 1

$ws.Range("C32").Value = 9
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 0.12

$ws.Range("C33").Value = 9
$ws.Range("D33").Value = 9
$ws.Range("E33").Value = 0.13

$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 0.06
